$d = $word.ActiveDocument

$d.Content.Find.Execute("167×7=", $true, $false, $false, $false, $false, $true, 1, $false, "473×9=", 2) | Out-Null
$d.Content.Find.Execute("850×2=", $true, $false, $false, $false, $false, $true, 1, $false, "944×7=", 2) | Out-Null
$d.Content.Find.Execute("978×4=", $true, $false, $false, $false, $false, $true, 1, $false, "251×7=", 2) | Out-Null
$d.Content.Find.Execute("370×8=", $true, $false, $false, $false, $false, $true, 1, $false, "556×6=", 2) | Out-Null
$d.Content.Find.Execute("357×3=", $true, $false, $false, $false, $false, $true, 1, $false, "432×5=", 2) | Out-Null
$d.Content.Find.Execute("836×6=", $true, $false, $false, $false, $false, $true, 1, $false, "283×8=", 2) | Out-Null
$d.Content.Find.Execute("285×3=", $true, $false, $false, $false, $false, $true, 1, $false, "601×7=", 2) | Out-Null
$d.Content.Find.Execute("476×6=", $true, $false, $false, $false, $false, $true, 1, $false, "418×8=", 2) | Out-Null
$d.Content.Find.Execute("399×6=", $true, $false, $false, $false, $false, $true, 1, $false, "788×5=", 2) | Out-Null
$d.Content.Find.Execute("701×4=", $true, $false, $false, $false, $false, $true, 1, $false, "638×2=", 2) | Out-Null
$d.Content.Find.Execute("242×7=", $true, $false, $false, $false, $false, $true, 1, $false, "511×7=", 2) | Out-Null
$d.Content.Find.Execute("634×6=", $true, $false, $false, $false, $false, $true, 1, $false, "353×7=", 2) | Out-Null
$d.Content.Find.Execute("201×6=", $true, $false, $false, $false, $false, $true, 1, $false, "323×3=", 2) | Out-Null
$d.Content.Find.Execute("119×5=", $true, $false, $false, $false, $false, $true, 1, $false, "598×3=", 2) | Out-Null
$d.Content.Find.Execute("209×2=", $true, $false, $false, $false, $false, $true, 1, $false, "647×8=", 2) | Out-Null
$d.Content.Find.Execute("815×2=", $true, $false, $false, $false, $false, $true, 1, $false, "531×5=", 2) | Out-Null
$d.Content.Find.Execute("258×5=", $true, $false, $false, $false, $false, $true, 1, $false, "438×5=", 2) | Out-Null
$d.Content.Find.Execute("237×2=", $true, $false, $false, $false, $false, $true, 1, $false, "146×6=", 2) | Out-Null
$d.Content.Find.Execute("814×7=", $true, $false, $false, $false, $false, $true, 1, $false, "366×4=", 2) | Out-Null
$d.Content.Find.Execute("868×4=", $true, $false, $false, $false, $false, $true, 1, $false, "126×2=", 2) | Out-Null
$d.Content.Find.Execute("237×4=", $true, $false, $false, $false, $false, $true, 1, $false, "589×3=", 2) | Out-Null
$d.Content.Find.Execute("587×8=", $true, $false, $false, $false, $false, $true, 1, $false, "597×9=", 2) | Out-Null
$d.Content.Find.Execute("427×5=", $true, $false, $false, $false, $false, $true, 1, $false, "267×8=", 2) | Out-Null
$d.Content.Find.Execute("991×3=", $true, $false, $false, $false, $false, $true, 1, $false, "432×7=", 2) | Out-Null
$d.Content.Find.Execute("598×9=", $true, $false, $false, $false, $false, $true, 1, $false, "582×3=", 2) | Out-Null
